$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.710.23'
$ws.Range('E2').Value = '  +2.20%  '
$ws.Range('D3').Value = '2.941.06'
$ws.Range('E3').Value = '  +0.43%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '592.50'
$ws.Range('E5').Value = '  -0.85%  '
$ws.Range('D6').Value = '147.39'
$ws.Range('E6').Value = '  +1.54%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = '2.939.86'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').Value = '0.505'
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +3.82%  '
$ws.Range('D11').Value = '0.148'
$ws.Range('E11').Value = '  +4.81%  '
$ws.Range('D12').Value = '0.440'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '0.0000234'
$ws.Range('E13').Value = '  +4.38%  '
$ws.Range('D14').Value = '32.58'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '3.427.63'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '62.707.28'
$ws.Range('E17').Value = '  +2.21%  '
$ws.Range('D18').Value = '6.67'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').Value = '2.940.81'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '439.30'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '13.34'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').Value = '0.664'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').Value = '7.02'
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '80.81'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').Value = '11.09'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').Value = '2.12'
$ws.Range('E26').Value = '  -2.80%  '
$ws.Range('D27').Value = '11.71'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = '7.15'
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('D31').Value = '2.60'
$ws.Range('E31').Value = '  -0.40%  '
$ws.Range('E32').Value = '  +14.13%  '
$ws.Range('D33').Value = '0.108'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = '26.29'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '0.990'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').Value = '3.07'
$ws.Range('E37').Value = '  +3.02%  '
$ws.Range('D38').Value = '5.55'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('D39').Value = '49.67'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').Value = '8.45'
$ws.Range('E41').Value = '  -1.38%  '
$ws.Range('E42').Value = '  -4.29%  '
$ws.Range('D43').Value = '0.279'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('D44').Value = '38.61'
$ws.Range('E44').Value = '  -8.75%  '
$ws.Range('D45').Value = '2.694.22'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '134.41'
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').Value = '359.72'
$ws.Range('E47').Value = '  -0.52%  '
$ws.Range('D48').Value = '0.0334'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '22.58'
$ws.Range('E51').Value = '  -4.23%  '
